# TestOverview.xlsx update:
#  - add "correct" markers (Comment / extra column) to a couple of existing
#    rows that were re-verified
#  - add two new test rows (LoadBalancing, UniReceive)
#  - widen the new "Comment" column (G) to fit its header
#  - leave the view scrolled down to the newly added rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- new column G width (header "Comment" needs more room) ---
$ws.Columns.Item(7).ColumnWidth = 10.5

# --- row 15 (Mutlicast / yes / no / no): mark as correct ---
$ws.Range("E15").Value = "correct"
$ws.Range("F15").Value = "correct"

# --- row 20 (Unicast / no / no / no): mark as correct ---
$ws.Range("E20").Value = "correct"
$ws.Range("F20").Value = "correct"

# --- row 22 (new): LoadBalancing test ---
$ws.Range("A22").Value = "LoadBalancing"
$ws.Range("B22").Value = "no"
$ws.Range("C22").Value = "no"
$ws.Range("D22").Value = "no"
$ws.Range("E22").Value = "correct"
$ws.Range("F22").Value = "correct"

# --- row 23 (new): UniReceive test ---
$ws.Range("A23").Value = "UniReceive"
$ws.Range("B23").Value = "no"
$ws.Range("C23").Value = "no"
$ws.Range("D23").Value = "no"

# --- scroll the view down and move the selection to the new row ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A22").Select()
